# API CONTRACT.docx change: "email" -> "id_user" for the
# /api/user/ verification row ("email", "token" body-request cell),
# matching the commit "change: API send verif / melakukan perubahan
# parameter pada API send verification email".
#
# Word's real Find/Replace in this runtime operates over the whole
# document even when invoked on a sub-Range, so instead of Find.Execute
# we locate the target table cell by its exact content and edit the
# Range directly - this keeps the edit scoped to the single cell that
# actually needs to change.

$d = $word.ActiveDocument

$targetCell = $null
$targetRow = 0
$targetCol = 0

for ($ti = 1; $ti -le $d.Tables.Count; $ti++) {
    $tbl = $d.Tables.Item($ti)
    for ($ri = 1; $ri -le $tbl.Rows.Count; $ri++) {
        $row = $tbl.Rows.Item($ri)
        for ($ci = 1; $ci -le $row.Cells.Count; $ci++) {
            $cell = $row.Cells.Item($ci)
            $cellText = $cell.Range.Text.TrimEnd([char]13, [char]7)
            if ($cellText -eq "“email”, “token”") {
                $targetCell = $cell
                $targetRow = $ri
                $targetCol = $ci
            }
        }
    }
}

if ($targetCell -eq $null) {
    Write-Output "ERROR: target cell not found"
} else {
    Write-Output "Found target cell: row $targetRow col $targetCol, text=[$($targetCell.Range.Text.TrimEnd([char]13,[char]7))]"

    # Locate the "email" word inside the cell (word 1 is the opening
    # curly quote, word 2 is "email").
    $words = $targetCell.Range.Words
    $emailWord = $words.Item(2)
    $startPos = $emailWord.Start

    if ($emailWord.Text -ne "email") {
        Write-Output "ERROR: expected word 'email', got '$($emailWord.Text)'"
    } else {
        $newWord = "id_user"

        # Replace "email" with "id_user" in place.
        $emailWord.Text = $newWord
        $endPos = $startPos + $newWord.Length

        # The runtime merges adjacent runs that end up with identical
        # formatting, which would leave "id_user" fused together with
        # the surrounding curly quotes in one run. Toggling a
        # (no-op) direct-formatting change on a sub-range forces the
        # run to be split at that range's boundaries while leaving the
        # visible formatting untouched, reproducing the three separate
        # runs ("“", "id_user", "”") the edit needs.

        # Split "id_user" away from the opening curly quote before it
        # and the closing curly quote after it.
        $innerRng = $d.Range($startPos, $endPos)
        $innerRng.Bold = 1
        $innerRng.Bold = 0

        # Split the closing curly quote away from the following
        # ", “token”" run so it becomes its own run too.
        $closingQuoteRng = $d.Range($endPos, $endPos + 1)
        $closingQuoteRng.Bold = 1
        $closingQuoteRng.Bold = 0

        Write-Output "After edit: [$($targetCell.Range.Text.TrimEnd([char]13,[char]7))]"
    }
}
